# SAO_PEDRO_DO_SUL.xlsx update:
#  - rename "Paineis DARQ" -> "PAINEIS DARQ"
#  - rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
#  - remove the "Desarquivamentos Pendentes" sheet entirely

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete() | Out-Null
$excel.DisplayAlerts = $true

# Keep the originally-selected tab ("PAINEIS DARQ") active/selected.
$wb.Worksheets.Item("PAINEIS DARQ").Activate()
